# Updated notebook, reran simulation
#
# The underlying simulation was rerun, which produced two brand-new HKL rows
# (inserted right after row 3 / "Spiral5") and shifted all the previously
# existing rows down by two. Along the way two new reflection labels
# ("Holden", "Rizzie Spiral") were introduced and "Thomas Hex" was renamed to
# "Matthies Hex".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# --- 1. Shift existing data rows 4..29 down to 6..31 (bottom-up so we never
#        clobber a row before it has been read). Row 9's label ("Thomas Hex")
#        is rewritten to "Matthies Hex" as it lands on its new row (11).
for ($r = 29; $r -ge 4; $r--) {
    $dst = $r + 2
    foreach ($col in $cols) {
        $srcAddr = $col + $r
        $dstAddr = $col + $dst
        $val = $ws.Range($srcAddr).Value2
        if (($col -eq "B") -and ($r -eq 9) -and ($val -eq "Thomas Hex")) {
            $val = "Matthies Hex"
        }
        $ws.Range($dstAddr).Value2 = $val
    }
}

# --- 2. Give the two freshly-vacated rows (4 and 5) the same bold/centered/
#        bordered formatting as the rest of the "A" index column, by copying
#        the already-formatted A3 cell onto them (keeps the shared style
#        instead of synthesizing a brand-new one).
$ws.Range("A3").Copy($ws.Range("A4"))
$ws.Range("A3").Copy($ws.Range("A5"))

# --- 3. Populate the two new rows with the freshly re-simulated data.
$ws.Range("A4").Value2 = 2
$ws.Range("B4").Value2 = "Holden"
$ws.Range("C4").Value2 = 1.05049211298127
$ws.Range("D4").Value2 = 0.8910832627642765
$ws.Range("E4").Value2 = 1.442126569225874
$ws.Range("F4").Value2 = 0.8433516243493203
$ws.Range("G4").Value2 = 0.8910832627642765
$ws.Range("H4").Value2 = 0.8078462762188927
$ws.Range("I4").Value2 = 1.157765799863202
$ws.Range("J4").Value2 = 0.8078462762188927
$ws.Range("K4").Value2 = 0.8910832627642765
$ws.Range("L4").Value2 = 1.157765799863202
$ws.Range("M4").Value2 = 0.9828060380410474
$ws.Range("N4").Value2 = 0.9828060380410474
$ws.Range("O4").Value2 = 0.9363212334771384
$ws.Range("P4").Value2 = 0.952231779615457
$ws.Range("Q4").Value2 = 0.9522317796154572
$ws.Range("R4").Value2 = 0.9369446504026621
$ws.Range("S4").Value2 = 0.9369446504026621
$ws.Range("T4").Value2 = 1.032110940900473

$ws.Range("A5").Value2 = 3
$ws.Range("B5").Value2 = "Rizzie Spiral"
$ws.Range("C5").Value2 = 1.361980976844551
$ws.Range("D5").Value2 = 2.334987030732487
$ws.Range("E5").Value2 = 0.4248095724664603
$ws.Range("F5").Value2 = 0.2523059868903106
$ws.Range("G5").Value2 = 2.334987030732487
$ws.Range("H5").Value2 = 0.004985389454532943
$ws.Range("I5").Value2 = 1.089124445491166
$ws.Range("J5").Value2 = 0.004985389454532943
$ws.Range("K5").Value2 = 2.334987030732487
$ws.Range("L5").Value2 = 1.089124445491166
$ws.Range("M5").Value2 = 0.5470549174728496
$ws.Range("N5").Value2 = 0.5470549174728496
$ws.Range("O5").Value2 = 0.4488052739453366
$ws.Range("P5").Value2 = 1.143032288559395
$ws.Range("Q5").Value2 = 1.143032288559396
$ws.Range("R5").Value2 = 1.441020974102668
$ws.Range("S5").Value2 = 1.441020974102668
$ws.Range("T5").Value2 = 0.9113655669799181
